# Adds a new "Inclusion?" column (K) ahead of the existing "Inclusion comments"
# column (shifted from K to L), and fills in inclusion-review comments for the
# hemi patients.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing "Inclusion comments" column (K),
# pushing it (and its data) one column to the right, into L.
$ws.Columns("K").Insert()

# --- Header row ---
$ws.Range("K1").Value = "Inclusion?"
$ws.Range("K1").HorizontalAlignment = -4108
$ws.Range("K1").Font.Bold = $true

$ws.Range("L1").Value = "Inclusion comments"
$ws.Range("L1").HorizontalAlignment = -4108

# --- Row 2 (12519) ---
$ws.Range("L2").Value = "fluid build up instead of brain expansion"
$ws.Range("L2").HorizontalAlignment = -4108

# --- Row 3 (13198) ---
$ws.Range("L3").Value = "ventricle enlargement (ventriculomegaly // hydrocephalus)"

# --- Row 6 (13782) ---
$ws.Range("L6").Value = "ventricle enlargement (ventriculomegaly // hydrocephalus)"

# --- Row 7 (13990) ---
$ws.Range("K7").Value = "?"
$ws.Range("K7").HorizontalAlignment = -4108
$ws.Range("L7").Value = "acute yes, fast ?? Skull edges difficult (skull edges always difficult)"

# --- Row 8 (14324) ---
$ws.Range("K8").Value = "?"
$ws.Range("K8").HorizontalAlignment = -4108
$ws.Range("L8").Value = "Skull can easily be seen, non symmetrical expansion"

# --- Row 10 (16754) ---
$ws.Range("K10").Value = "?"
$ws.Range("K10").HorizontalAlignment = -4108
$ws.Range("L10").Value = "Interesting compression - the skin follows shape of symmetrical expansion but lesion is pressing on expansion as well. "

# --- Row 12 (19344) ---
$ws.Range("K12").Value = 1
$ws.Range("K12").HorizontalAlignment = -4108
$ws.Range("L12").Value = "Skull visible, mostly symmetrical expansion considering nothing is spherical"

# --- Row 13 (19575) ---
$ws.Range("L13").Value = "Not much expansion"

# --- Row 14 (19978) ---
$ws.Range("K14").Value = 1
$ws.Range("K14").HorizontalAlignment = -4108
$ws.Range("L14").Value = "hard to see skull anchor points but makes a very nice shape, little else confounding the image"

# --- Row 15 (19981) ---
$ws.Range("L15").Value = "like 12519, fluid build up instead of brain expansion."

# --- Row 16 (20174) ---
$ws.Range("K16").Value = 1
$ws.Range("K16").HorizontalAlignment = -4108
$ws.Range("L16").Value = "ultra-fast to fast only - acute is wildly deformed. 3 months is okay. "

# --- Row 17 (20651) ---
$ws.Range("L17").Value = "lesion layer on top of free bulge"

# --- Row 19 (20942) ---
$ws.Range("K19").Value = "?"
$ws.Range("K19").HorizontalAlignment = -4108
$ws.Range("L19").Value = "non symmetrical expansion, easy to see skull points"

# --- Row 20 (21221) ---
$ws.Range("L20").Value = "brain didn’t change shape sufficiently"

# --- Row 23 (22725) ---
$ws.Range("K23").Value = 1
$ws.Range("K23").HorizontalAlignment = -4108
$ws.Range("L23").Value = "!! THE BEST ONE"

# --- Row 24 (22785) --- (comment typed before the "??tech" inclusion flag)
$ws.Range("L24").Value = "could be good - bad registration"
$ws.Range("K24").Value = "??tech"
$ws.Range("K24").HorizontalAlignment = -4108

# --- Row 25 (23348) ---
$ws.Range("L25").Value = "push back on expansion in sigmoid shape -> fills out over time. 3 months = ventriculomegaly "

# Match the author's final selection position.
$ws.Range("L26").Select()
